# Auto-generated Excel COM-interop script to update TPM values
# as described by the commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.05089966666666667
$ws.Range("H2").Value = 0.152699
$ws.Range("I2").Value = 0.02671091810242436
$ws.Range("J2").Value = 0.03728162213961778
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6435283333333334
$ws.Range("N2").Value = 1.930585
$ws.Range("O2").Value = 0.1384187503011309
$ws.Range("P2").Value = 0.1430362385488992
$ws.Range("Q2").Value = 0.03275537765722222
$ws.Range("R2").Value = 0.294798398915
$ws.Range("S2").Value = 0.003697291903133434
$ws.Range("T2").Value = 0.005332622997852291
$ws.Range("G3").Value = 0.05089966666666667
$ws.Range("H3").Value = 0.152699
$ws.Range("I3").Value = 0.02671091810242436
$ws.Range("J3").Value = 0.03728162213961778
$ws.Range("O3").Value = 0.2324275794584003
$ws.Range("P3").Value = 0.2401810927235583
$ws.Range("Q3").Value = 0.05500160293711111
$ws.Range("R3").Value = 0.495014426434
$ws.Range("S3").Value = 0.006208354039658063
$ws.Range("T3").Value = 0.008954340744000202
$ws.Range("G4").Value = 0.05089966666666667
$ws.Range("H4").Value = 0.152699
$ws.Range("I4").Value = 0.02671091810242436
$ws.Range("J4").Value = 0.03728162213961778
$ws.Range("M4").Value = 1.770781666666667
$ws.Range("N4").Value = 5.312345000000001
$ws.Range("O4").Value = 0.3808835954223518
$ws.Range("P4").Value = 0.3935894284240538
$ws.Range("Q4").Value = 0.0901321965727778
$ws.Range("R4").Value = 0.8111897691550001
$ws.Range("S4").Value = 0.01017375052388337
$ws.Range("T4").Value = 0.01467365234865371
$ws.Range("G5").Value = 0.05089966666666667
$ws.Range("H5").Value = 0.152699
$ws.Range("I5").Value = 0.02671091810242436
$ws.Range("J5").Value = 0.03728162213961778
$ws.Range("M5").Value = 0.45025
$ws.Range("N5").Value = 0.9005000000000001
$ws.Range("O5").Value = 0.09684584049355638
$ws.Range("P5").Value = 0.06671766993594362
$ws.Range("Q5").Value = 0.02291757491666667
$ws.Range("R5").Value = 0.1375054495
$ws.Range("S5").Value = 0.002586841313983838
$ws.Range("T5").Value = 0.002487342960587587
$ws.Range("G6").Value = 0.05089966666666667
$ws.Range("H6").Value = 0.152699
$ws.Range("I6").Value = 0.02671091810242436
$ws.Range("J6").Value = 0.03728162213961778
$ws.Range("M6").Value = 0.7039926666666667
$ws.Range("N6").Value = 2.111978
$ws.Range("O6").Value = 0.1514242343245606
$ws.Range("P6").Value = 0.1564755703675451
$ws.Range("Q6").Value = 0.03583299206911111
$ws.Range("R6").Value = 0.322496928622
$ws.Range("S6").Value = 0.004044680321765654
$ws.Range("T6").Value = 0.00583366308852399
$ws.Range("I7").Value = 0.1226793098007496
$ws.Range("J7").Value = 0.1712289953794413
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6435283333333334
$ws.Range("N7").Value = 1.930585
$ws.Range("O7").Value = 0.1384187503011309
$ws.Range("P7").Value = 0.1430362385488992
$ws.Range("Q7").Value = 0.1504406216155556
$ws.Range("R7").Value = 1.35396559454
$ws.Range("S7").Value = 0.01698111675042504
$ws.Range("T7").Value = 0.02449195142958212
$ws.Range("I8").Value = 0.1226793098007496
$ws.Range("J8").Value = 0.1712289953794413
$ws.Range("O8").Value = 0.2324275794584003
$ws.Range("P8").Value = 0.2401810927235583
$ws.Range("S8").Value = 0.02851405502661544
$ws.Range("T8").Value = 0.04112596721619132
$ws.Range("I9").Value = 0.1226793098007496
$ws.Range("J9").Value = 0.1712289953794413
$ws.Range("M9").Value = 1.770781666666667
$ws.Range("N9").Value = 5.312345000000001
$ws.Range("O9").Value = 0.3808835954223518
$ws.Range("P9").Value = 0.3935894284240538
$ws.Range("Q9").Value = 0.4139638938644445
$ws.Range("R9").Value = 3.725675044780001
$ws.Range("S9").Value = 0.04672653660084208
$ws.Range("T9").Value = 0.06739392242101924
$ws.Range("I10").Value = 0.1226793098007496
$ws.Range("J10").Value = 0.1712289953794413
$ws.Range("M10").Value = 0.45025
$ws.Range("N10").Value = 0.9005000000000001
$ws.Range("O10").Value = 0.09684584049355638
$ws.Range("P10").Value = 0.06671766993594362
$ws.Range("Q10").Value = 0.1052570436666667
$ws.Range("R10").Value = 0.6315422620000001
$ws.Range("S10").Value = 0.01188098086882298
$ws.Range("T10").Value = 0.01142399959718878
$ws.Range("I11").Value = 0.1226793098007496
$ws.Range("J11").Value = 0.1712289953794413
$ws.Range("M11").Value = 0.7039926666666667
$ws.Range("N11").Value = 2.111978
$ws.Range("O11").Value = 0.1514242343245606
$ws.Range("P11").Value = 0.1564755703675451
$ws.Range("Q11").Value = 0.1645756509857778
$ws.Range("R11").Value = 1.481180858872
$ws.Range("S11").Value = 0.01857662055404407
$ws.Range("T11").Value = 0.02679315471545982
$ws.Range("G12").Value = 1.620901
$ws.Range("H12").Value = 3.241802
$ws.Range("I12").Value = 0.8506097720968261
$ws.Range("J12").Value = 0.791489382480941
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.6435283333333334
$ws.Range("N12").Value = 1.930585
$ws.Range("O12").Value = 0.1384187503011309
$ws.Range("P12").Value = 0.1430362385488992
$ws.Range("Q12").Value = 1.043095719028333
$ws.Range("R12").Value = 6.25857431417
$ws.Range("S12").Value = 0.1177403416475724
$ws.Range("T12").Value = 0.1132116641214648
$ws.Range("G13").Value = 1.620901
$ws.Range("H13").Value = 3.241802
$ws.Range("I13").Value = 0.8506097720968261
$ws.Range("J13").Value = 0.791489382480941
$ws.Range("O13").Value = 0.2324275794584003
$ws.Range("P13").Value = 0.2401810927235583
$ws.Range("Q13").Value = 1.751527250388667
$ws.Range("R13").Value = 10.509163502332
$ws.Range("S13").Value = 0.1977051703921268
$ws.Range("T13").Value = 0.1901007847633668
$ws.Range("G14").Value = 1.620901
$ws.Range("H14").Value = 3.241802
$ws.Range("I14").Value = 0.8506097720968261
$ws.Range("J14").Value = 0.791489382480941
$ws.Range("M14").Value = 1.770781666666667
$ws.Range("N14").Value = 5.312345000000001
$ws.Range("O14").Value = 0.3808835954223518
$ws.Range("P14").Value = 0.3935894284240538
$ws.Range("Q14").Value = 2.870261774281667
$ws.Range("R14").Value = 17.22157064569
$ws.Range("S14").Value = 0.3239833082976264
$ws.Range("T14").Value = 0.3115218536543809
$ws.Range("G15").Value = 1.620901
$ws.Range("H15").Value = 3.241802
$ws.Range("I15").Value = 0.8506097720968261
$ws.Range("J15").Value = 0.791489382480941
$ws.Range("M15").Value = 0.45025
$ws.Range("N15").Value = 0.9005000000000001
$ws.Range("O15").Value = 0.09684584049355638
$ws.Range("P15").Value = 0.06671766993594362
$ws.Range("Q15").Value = 0.72981067525
$ws.Range("R15").Value = 2.919242701
$ws.Range("S15").Value = 0.08237801831074956
$ws.Range("T15").Value = 0.05280632737816726
$ws.Range("G16").Value = 1.620901
$ws.Range("H16").Value = 3.241802
$ws.Range("I16").Value = 0.8506097720968261
$ws.Range("J16").Value = 0.791489382480941
$ws.Range("M16").Value = 0.7039926666666667
$ws.Range("N16").Value = 2.111978
$ws.Range("O16").Value = 0.1514242343245606
$ws.Range("P16").Value = 0.1564755703675451
$ws.Range("Q16").Value = 1.141102417392667
$ws.Range("R16").Value = 6.846614504356
$ws.Range("S16").Value = 0.1288029334487509
$ws.Range("T16").Value = 0.1238487525635613
